$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(68, 8).Value = 69999  # H68
$ws.Cells.Item(68, 10).Value = 69999  # J68
$ws.Cells.Item(68, 12).Value = 69999  # L68
$ws.Cells.Item(68, 14).Value = -71497  # N68
$ws.Cells.Item(71, 8).Value = 69999  # H71
$ws.Cells.Item(71, 10).Value = 69999  # J71
$ws.Cells.Item(71, 12).Value = 209997  # L71
$ws.Cells.Item(71, 14).Value = -217485  # N71
$ws.Cells.Item(76, 8).Value = 3693.2856  # H76
$ws.Cells.Item(76, 9).Value = 4000  # I76
$ws.Cells.Item(76, 10).Value = 3570.6  # J76
$ws.Cells.Item(76, 11).Value = 4000  # K76
$ws.Cells.Item(76, 12).Value = 3570.6  # L76
$ws.Cells.Item(76, 13).Value = -3685  # M76
$ws.Cells.Item(76, 14).Value = -4200.6  # N76
$ws.Cells.Item(79, 8).Value = 3693.2856  # H79
$ws.Cells.Item(79, 9).Value = 4000  # I79
$ws.Cells.Item(79, 10).Value = 3570.6  # J79
$ws.Cells.Item(79, 11).Value = 4000  # K79
$ws.Cells.Item(79, 12).Value = 3570.6  # L79
$ws.Cells.Item(79, 13).Value = -2908  # M79
$ws.Cells.Item(79, 14).Value = -5754.6  # N79
$ws.Cells.Item(98, 8).Value = 1120  # H98
$ws.Cells.Item(98, 9).Value = 1049.1111  # I98
$ws.Cells.Item(98, 10).Value = 1332.6666  # J98
$ws.Cells.Item(98, 11).Value = 1049.1111  # K98
$ws.Cells.Item(98, 12).Value = 1332.6666  # L98
$ws.Cells.Item(98, 13).Value = 448.8888999999999  # M98
$ws.Cells.Item(98, 14).Value = -4328.6666  # N98
$ws.Cells.Item(100, 8).Value = 1174.625  # H100
$ws.Cells.Item(100, 9).Value = 651.2  # I100
$ws.Cells.Item(100, 10).Value = 2047  # J100
$ws.Cells.Item(100, 11).Value = 651.2  # K100
$ws.Cells.Item(100, 12).Value = 2047  # L100
$ws.Cells.Item(100, 13).Value = -110.2  # M100
$ws.Cells.Item(100, 14).Value = -3129  # N100
$ws.Cells.Item(111, 8).Value = 80277.30499999999  # H111
$ws.Cells.Item(111, 9).Value = 1529.625  # I111
$ws.Cells.Item(111, 11).Value = 4588.875  # K111
$ws.Cells.Item(111, 13).Value = -1521.875  # M111
$ws.Cells.Item(122, 8).Value = 1120  # H122
$ws.Cells.Item(122, 9).Value = 1049.1111  # I122
$ws.Cells.Item(122, 10).Value = 1332.6666  # J122
$ws.Cells.Item(122, 11).Value = 3147.3333  # K122
$ws.Cells.Item(122, 12).Value = 3997.9998  # L122
$ws.Cells.Item(122, 13).Value = -697.3333000000002  # M122
$ws.Cells.Item(122, 14).Value = -8897.9998  # N122
$ws.Cells.Item(141, 8).Value = 5277.4443  # H141
$ws.Cells.Item(141, 9).Value = 3845.6924  # I141
$ws.Cells.Item(141, 10).Value = 9000  # J141
$ws.Cells.Item(141, 11).Value = 11537.0772  # K141
$ws.Cells.Item(141, 12).Value = 27000  # L141
$ws.Cells.Item(141, 13).Value = -6357.0772  # M141
$ws.Cells.Item(141, 14).Value = -37360  # N141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7181.3027  # H32
$ws.Cells.Item(32, 9).Value = 2939.7856  # I32
$ws.Cells.Item(32, 11).Value = 2939.7856  # K32
$ws.Cells.Item(32, 13).Value = -2652.7856  # M32
$ws.Cells.Item(102, 8).Value = 0  # H102
$ws.Cells.Item(102, 9).Value = 0  # I102
$ws.Cells.Item(102, 11).Value = 0  # K102
$ws.Cells.Item(102, 13).ClearContents()  # M102
$ws.Cells.Item(110, 8).Value = 2783.6  # H110
$ws.Cells.Item(110, 10).Value = 964.5  # J110
$ws.Cells.Item(110, 12).Value = 964.5  # L110
$ws.Cells.Item(110, 14).Value = -5054.5  # N110
$ws.Cells.Item(132, 8).Value = 3455.2222  # H132
$ws.Cells.Item(132, 9).Value = 1819.3611  # I132
$ws.Cells.Item(132, 11).Value = 5458.0833  # K132
$ws.Cells.Item(132, 13).Value = -2928.0833  # M132
$ws.Cells.Item(139, 8).Value = 112723  # H139
$ws.Cells.Item(139, 10).Value = 112723  # J139
$ws.Cells.Item(139, 12).Value = 112723  # L139
$ws.Cells.Item(139, 14).Value = -123003  # N139

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 742.6429000000001  # H94
$ws.Cells.Item(94, 9).Value = 672.63635  # I94
$ws.Cells.Item(94, 11).Value = 672.63635  # K94
$ws.Cells.Item(94, 13).Value = -221.63635  # M94
$ws.Cells.Item(105, 8).Value = 6894.8  # H105
$ws.Cells.Item(105, 10).Value = 6666.6665  # J105
$ws.Cells.Item(105, 12).Value = 6666.6665  # L105
$ws.Cells.Item(105, 14).Value = -10160.6665  # N105
$ws.Cells.Item(107, 8).Value = 2033.3793  # H107
$ws.Cells.Item(107, 9).Value = 1590.8  # I107
$ws.Cells.Item(107, 10).Value = 4799.5  # J107
$ws.Cells.Item(107, 11).Value = 1590.8  # K107
$ws.Cells.Item(107, 12).Value = 4799.5  # L107
$ws.Cells.Item(107, 13).Value = 329.2  # M107
$ws.Cells.Item(107, 14).Value = -8639.5  # N107

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 736.8333  # H16
$ws.Cells.Item(16, 9).Value = 584.2  # I16
$ws.Cells.Item(16, 11).Value = 584.2  # K16
$ws.Cells.Item(16, 13).Value = -297.2  # M16
$ws.Cells.Item(22, 8).Value = 1977.5625  # H22
$ws.Cells.Item(22, 9).Value = 1723.4445  # I22
$ws.Cells.Item(22, 11).Value = 1723.4445  # K22
$ws.Cells.Item(22, 13).Value = -1373.4445  # M22
$ws.Cells.Item(31, 8).Value = 7962.6978  # H31
$ws.Cells.Item(31, 9).Value = 2793.1428  # I31
$ws.Cells.Item(31, 11).Value = 2793.1428  # K31
$ws.Cells.Item(31, 13).Value = -2498.1428  # M31
$ws.Cells.Item(34, 8).Value = 7962.6978  # H34
$ws.Cells.Item(34, 9).Value = 2793.1428  # I34
$ws.Cells.Item(34, 11).Value = 2793.1428  # K34
$ws.Cells.Item(34, 13).Value = -2591.1428  # M34
$ws.Cells.Item(95, 8).Value = 42424.8  # H95
$ws.Cells.Item(95, 10).Value = 42424.8  # J95
$ws.Cells.Item(95, 12).Value = 42424.8  # L95
$ws.Cells.Item(95, 14).Value = -47916.8  # N95
$ws.Cells.Item(99, 8).Value = 3729.8572  # H99
$ws.Cells.Item(99, 9).Value = 4018.1667  # I99
$ws.Cells.Item(99, 11).Value = 4018.1667  # K99
$ws.Cells.Item(99, 13).Value = -2520.1667  # M99
$ws.Cells.Item(113, 8).Value = 736.8333  # H113
$ws.Cells.Item(113, 9).Value = 584.2  # I113
$ws.Cells.Item(113, 11).Value = 584.2  # K113
$ws.Cells.Item(113, 13).Value = 1585.8  # M113
$ws.Cells.Item(126, 8).Value = 3729.8572  # H126
$ws.Cells.Item(126, 9).Value = 4018.1667  # I126
$ws.Cells.Item(126, 11).Value = 12054.5001  # K126
$ws.Cells.Item(126, 13).Value = -9584.500100000001  # M126
$ws.Cells.Item(134, 8).Value = 8176.6816  # H134
$ws.Cells.Item(134, 9).Value = 7210.4443  # I134
$ws.Cells.Item(134, 10).Value = 8845.615  # J134
$ws.Cells.Item(134, 11).Value = 21631.3329  # K134
$ws.Cells.Item(134, 12).Value = 26536.845  # L134
$ws.Cells.Item(134, 13).Value = -19096.3329  # M134
$ws.Cells.Item(134, 14).Value = -31606.845  # N134

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 70.26316  # H38
$ws.Cells.Item(38, 10).Value = 85.5  # J38
$ws.Cells.Item(38, 12).Value = 256.5  # L38
$ws.Cells.Item(38, 14).Value = -950.5  # N38
$ws.Cells.Item(51, 8).Value = 799.5  # H51
$ws.Cells.Item(51, 10).Value = 99  # J51
$ws.Cells.Item(51, 12).Value = 297  # L51
$ws.Cells.Item(51, 14).Value = -1217  # N51
$ws.Cells.Item(64, 8).Value = 8550  # H64
$ws.Cells.Item(64, 9).Value = 8550  # I64
$ws.Cells.Item(64, 11).Value = 25650  # K64
$ws.Cells.Item(64, 13).Value = -25380  # M64
$ws.Cells.Item(67, 8).Value = 8550  # H67
$ws.Cells.Item(67, 9).Value = 8550  # I67
$ws.Cells.Item(67, 11).Value = 25650  # K67
$ws.Cells.Item(67, 13).Value = -24714  # M67
$ws.Cells.Item(97, 8).Value = 1122.75  # H97
$ws.Cells.Item(97, 9).Value = 1010.5  # I97
$ws.Cells.Item(97, 11).Value = 3031.5  # K97
$ws.Cells.Item(97, 13).Value = -2535.5  # M97
$ws.Cells.Item(133, 8).Value = 5811.4614  # H133
$ws.Cells.Item(133, 9).Value = 3913.5454  # I133
$ws.Cells.Item(133, 10).Value = 16250  # J133
$ws.Cells.Item(133, 11).Value = 11740.6362  # K133
$ws.Cells.Item(133, 12).Value = 48750  # L133
$ws.Cells.Item(133, 13).Value = -6680.636200000001  # M133
$ws.Cells.Item(133, 14).Value = -58870  # N133
$ws.Cells.Item(134, 8).Value = 1759.875  # H134
$ws.Cells.Item(134, 9).Value = 1759.875  # I134
$ws.Cells.Item(134, 11).Value = 5279.625  # K134
$ws.Cells.Item(134, 13).Value = -209.625  # M134
$ws.Cells.Item(140, 8).Value = 4664.476  # H140
$ws.Cells.Item(140, 9).Value = 27640  # I140
$ws.Cells.Item(140, 11).Value = 82920  # K140
$ws.Cells.Item(140, 13).Value = -77740  # M140
$ws.Cells.Item(141, 8).Value = 6038.8945  # H141
$ws.Cells.Item(141, 9).Value = 4969.5713  # I141
$ws.Cells.Item(141, 11).Value = 14908.7139  # K141
$ws.Cells.Item(141, 13).Value = -9728.713899999999  # M141

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3190.1765  # H80
$ws.Cells.Item(80, 9).Value = 3148  # I80
$ws.Cells.Item(80, 10).Value = 3237.625  # J80
$ws.Cells.Item(80, 11).Value = 3148  # K80
$ws.Cells.Item(80, 12).Value = 3237.625  # L80
$ws.Cells.Item(80, 13).Value = -2150  # M80
$ws.Cells.Item(80, 14).Value = -5233.625  # N80
$ws.Cells.Item(83, 8).Value = 3190.1765  # H83
$ws.Cells.Item(83, 9).Value = 3148  # I83
$ws.Cells.Item(83, 10).Value = 3237.625  # J83
$ws.Cells.Item(83, 11).Value = 15740  # K83
$ws.Cells.Item(83, 12).Value = 16188.125  # L83
$ws.Cells.Item(83, 13).Value = -10748  # M83
$ws.Cells.Item(83, 14).Value = -26172.125  # N83

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 1747.5161  # H93
$ws.Cells.Item(93, 9).Value = 1585.0526  # I93
$ws.Cells.Item(93, 10).Value = 2004.75  # J93
$ws.Cells.Item(93, 11).Value = 1585.0526  # K93
$ws.Cells.Item(93, 12).Value = 2004.75  # L93
$ws.Cells.Item(93, 13).Value = -337.0526  # M93
$ws.Cells.Item(93, 14).Value = -4500.75  # N93
$ws.Cells.Item(136, 8).Value = 8240.421  # H136
$ws.Cells.Item(136, 9).Value = 4273.2173  # I136
$ws.Cells.Item(136, 10).Value = 9962.038  # J136
$ws.Cells.Item(136, 11).Value = 12819.6519  # K136
$ws.Cells.Item(136, 12).Value = 29886.114  # L136
$ws.Cells.Item(136, 13).Value = -10269.6519  # M136
$ws.Cells.Item(136, 14).Value = -34986.114  # N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 37494.484  # H81
$ws.Cells.Item(81, 9).Value = 68832.2  # I81
$ws.Cells.Item(81, 11).Value = 137664.4  # K81
$ws.Cells.Item(81, 13).Value = -136603.4  # M81
$ws.Cells.Item(84, 8).Value = 37494.484  # H84
$ws.Cells.Item(84, 9).Value = 68832.2  # I84
$ws.Cells.Item(84, 11).Value = 688322  # K84
$ws.Cells.Item(84, 13).Value = -683018  # M84
$ws.Cells.Item(107, 8).Value = 1690.9286  # H107
$ws.Cells.Item(107, 9).Value = 1725.4  # I107
$ws.Cells.Item(107, 11).Value = 5176.200000000001  # K107
$ws.Cells.Item(107, 13).Value = -3256.200000000001  # M107
$ws.Cells.Item(113, 8).Value = 855.7059  # H113
$ws.Cells.Item(113, 9).Value = 962.8461  # I113
$ws.Cells.Item(113, 10).Value = 507.5  # J113
$ws.Cells.Item(113, 11).Value = 2888.5383  # K113
$ws.Cells.Item(113, 12).Value = 1522.5  # L113
$ws.Cells.Item(113, 13).Value = -718.5383000000002  # M113
$ws.Cells.Item(113, 14).Value = -5862.5  # N113
